$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Full target header row (A1:AE1) ---
# Two new reactant columns (R7, R8) are inserted before W1, and two new
# species columns (X7, X8) are inserted before "Zeroes".
$headers = @("AC","NumReact","F1","F2","F3","M1","M2","M3","R1","R2","R3","R4","R5","R6","R7","R8","W1","W2","W3","EM1","EM2","EM3","X1","X2","X3","X4","X5","X6","X7","X8","Zeroes")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Re-apply the header style (bold font, thin border, centered/top aligned - same
# style as the original header cells) across the whole header row, since the
# newly created cells (AB1..AE1) would otherwise default to the plain style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1:AE1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- New data rows 2 and 3 (A..AE = 31 columns each) ---
$row2 = @(1,3,0,0,0,1,1,1,1,0,0,1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,3)
$row3 = @(2,3,0,0,0,1,1,1,1,0,0,1,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,3)

for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $row2[$i]
}
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# Column A on rows 2 and 3 uses the same (bold/border/center) style as the header.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Restore the numeric values for A2/A3 (PasteSpecial(-4122) only copies formats, but
# set them explicitly again to guarantee correctness).
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
